# Patton's Best - Events.xlsx edit
# "Got initial combat caledar check loop working"
#
# Semantic changes applied:
#  - Row 8 (event "e008") renamed/reworked to "e007": fixes the weather-roll
#    button/image names (tWeather -> Weather, DiceRoll -> DieRoll) and
#    re-wraps the paragraph so the button/image sit on their own lines.
#  - Row 7 (event "e006" Combat Calendar Check) text re-wrapped so the
#    DieRoll image sits on its own line, and the trailing blank line is
#    removed.
#  - Row 8 grows a line, so its row height increases to match row 4/5.
#  - View state (selection, scroll position) nudged per the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# --- Row 7: e006 Combat Calendar Check -> minor text re-wrap -------------
$ws.Range("B7").Value = @'
<Bold>e006 Combat Calendar Check</Bold> 
<InlineUIContainer><Button Content='r4.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<InlineUIContainer><Button Content='Calendar' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
Roll for possible combat today. If die &lt;= probability, start morning briefing per 
<InlineUIContainer><Button Content='e007' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.  
Otherwise continue with next day check.
<LineBreak/><LineBreak/>
Date from Combat Calendar: DATE<LineBreak/>
Expected Resistance: RESISTANCE<LineBreak/>
Probablility of Combat: PROBABILITY &gt;= 
<InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/>
'@

# --- Row 8: e008 -> e007 (Morning Briefing - Weather Roll) ----------------
$ws.Range("A8").Value = "e007"
$ws.Range("B8").Value = @'
<Bold>e007 Morning Briefing - Weather Roll</Bold> 
<LineBreak/><LineBreak/>
Consult weather table 
<InlineUIContainer><Button Content='Weather' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
after rolling to determine weather for today:
<InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/>
'@

# Row 8 picked up an extra wrapped line from the new text -> taller row.
$ws.Rows.Item(8).RowHeight = 99.85

# --- View state nudges (scroll position / selection) ----------------------
$ws.Range("E7").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$win.Left = 25974 | Out-Null
